$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.185.68"
$ws.Range("E2").Value = "  -0.51%  "

$ws.Range("D3").Value = "3.014.14"
$ws.Range("E3").Value = "  +0.14%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.54%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.22%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "3.014.58"
$ws.Range("E8").Value = "  +0.21%  "

$ws.Range("E9").Value = "  -1.60%  "

$ws.Range("E10").Value = "  +8.75%  "

$ws.Range("E11").Value = "  +0.48%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.459"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.05%  "

$ws.Range("E13").Value = "  +0.79%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.67%  "

$ws.Range("E15").Value = "  +2.47%  "

$ws.Range("D16").Value = "3.512.01"
$ws.Range("E16").Value = "  +0.02%  "

$ws.Range("D17").Value = "62.188.81"
$ws.Range("E17").Value = "  -0.44%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.59%  "

$ws.Range("D19").Value = "3.010.30"
$ws.Range("E19").Value = "  -0.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "446.49"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.81%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.56%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.689"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.08%  "

$ws.Range("E23").Value = "  -0.25%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.53%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +10.23%  "

$ws.Range("E26").Value = "  +1.11%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.19%  "

$ws.Range("E28").Value = "  +0.07%  "

$ws.Range("E29").Value = "  +1.97%  "

$ws.Range("E31").Value = "  +2.43%  "

$ws.Range("E32").Value = "  +0.99%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.46"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.37%  "

$ws.Range("E34").Value = "  +0.45%  "

$ws.Range("D35").Value = "0.0₃0850"
$ws.Range("E35").Value = "  +4.32%  "

$ws.Range("E36").Value = "  +0.21%  "

$ws.Range("E37").Value = "  +1.26%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "50.16"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.42%  "

$ws.Range("E39").Value = "  -0.88%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.05"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.42%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.96"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.67%  "

$ws.Range("E42").Value = "  -0.46%  "

$ws.Range("B43").Value = "Arweave"
$ws.Range("C43").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.25"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +10.25%  "

$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.285"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.08%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "394.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.79%  "

$ws.Range("E46").Value = "  -1.77%  "

$ws.Range("D47").Value = "2.726.63"
$ws.Range("E47").Value = "  -0.56%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.50"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.93%  "

$ws.Range("E49").Value = "  +0.09%  "

$ws.Range("E50").Value = "  -0.78%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.108"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.48%  "
